$d = $word.ActiveDocument

# 1) Update the date mentioned in the intro paragraph: 4 -> 7 de marzo de 2022
$d.Content.Find.Execute(
    "fecha  4 de marzo de 2022",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "fecha  7 de marzo de 2022", 2)

# 2) Add a new student row to the data table (2nd table in the document)
$t = $d.Tables(2)
$newRow = $t.Rows.Add()
$newRow.Cells(1).Range.Text = "Díez Viñas Malena"
$newRow.Cells(2).Range.Text = "14d"
$newRow.Cells(3).Range.Text = ""
$newRow.Cells(4).Range.Text = "edfsf"
$newRow.Cells(5).Range.Text = "400"
$newRow.Cells(6).Range.Text = "2022-03-08"
$newRow.Cells(7).Range.Text = "2022-03-09"

# 3) Update the signature date line: 4 -> 7 de marzo 2022
$d.Content.Find.Execute(
    "En Puertollano a  4  de marzo  2022",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "En Puertollano a  7  de marzo  2022", 2)
